# The "Date" column (BF) holds a malformed date string ("5-22-2013-14")
# on every data row. Fix it to the correct ISO-style date "2014-05-22"
# for each row, keeping the cell's type as plain text (not an Excel date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-22-2013-14"
$newValue = "2014-05-22"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1
$col = 58  # column BF

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq $oldValue) {
        # Force text formatting first so Excel doesn't reinterpret the
        # ISO-looking string as a date serial number, then strip the
        # number-format override back off so the cell's style is left
        # exactly as it was (default/general).
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    }
}
